$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns G:K hold numeric-looking text (run counts, strike rates, etc.)
# that must stay text (matching the rest of the sheet), so format as Text
# before writing the values -- otherwise Excel would coerce "50" etc. to numbers.
$ws.Range("G17:K31").NumberFormat = "@"

# Row 17
$ws.Range("A17").Value = ' Dubai (DSC)'
$ws.Range("B17").Value = ' October 25 2020'
$ws.Range("C17").Value = 'Super Kings won by 8 wickets (with 8 balls remaining)'
$ws.Range("D17").Value = 'Royal Challengers Bangalore'
$ws.Range("E17").Value = 'Chennai Super Kings'
$ws.Range("F17").Value = 'Virat Kohli (c)'
$ws.Range("G17").Value = '50'
$ws.Range("H17").Value = '43'
$ws.Range("I17").Value = '1'
$ws.Range("J17").Value = '1'
$ws.Range("K17").Value = '116.27'

# Row 18
$ws.Range("A18").Value = ' Abu Dhabi'
$ws.Range("B18").Value = ' October 28 2020'
$ws.Range("C18").Value = 'Mumbai won by 5 wickets (with 5 balls remaining)'
$ws.Range("D18").Value = 'Royal Challengers Bangalore'
$ws.Range("E18").Value = 'Mumbai Indians'
$ws.Range("F18").Value = 'Virat Kohli (c)'
$ws.Range("G18").Value = '9'
$ws.Range("H18").Value = '14'
$ws.Range("I18").Value = '0'
$ws.Range("J18").Value = '0'
$ws.Range("K18").Value = '64.28'

# Row 19
$ws.Range("A19").Value = ' Abu Dhabi'
$ws.Range("B19").Value = ' November 02 2020'
$ws.Range("C19").Value = 'Capitals won by 6 wickets (with 6 balls remaining)'
$ws.Range("D19").Value = 'Royal Challengers Bangalore'
$ws.Range("E19").Value = 'Delhi Capitals'
$ws.Range("F19").Value = 'Virat Kohli (c)'
$ws.Range("G19").Value = '29'
$ws.Range("H19").Value = '24'
$ws.Range("I19").Value = '2'
$ws.Range("J19").Value = '1'
$ws.Range("K19").Value = '120.83'

# Row 20
$ws.Range("A20").Value = ' Dubai (DSC)'
$ws.Range("B20").Value = ' September 24 2020'
$ws.Range("C20").Value = 'Kings XI won by 97 runs'
$ws.Range("D20").Value = 'Royal Challengers Bangalore'
$ws.Range("E20").Value = 'Kings XI Punjab'
$ws.Range("F20").Value = 'Virat Kohli (c)'
$ws.Range("G20").Value = '1'
$ws.Range("H20").Value = '5'
$ws.Range("I20").Value = '0'
$ws.Range("J20").Value = '0'
$ws.Range("K20").Value = '20.00'

# Row 21
$ws.Range("A21").Value = ' Abu Dhabi'
$ws.Range("B21").Value = ' November 06 2020'
$ws.Range("C21").Value = 'Sunrisers won by 6 wickets (with 2 balls remaining)'
$ws.Range("D21").Value = 'Royal Challengers Bangalore'
$ws.Range("E21").Value = 'Sunrisers Hyderabad'
$ws.Range("F21").Value = 'Virat Kohli (c)'
$ws.Range("G21").Value = '6'
$ws.Range("H21").Value = '7'
$ws.Range("I21").Value = '0'
$ws.Range("J21").Value = '0'
$ws.Range("K21").Value = '85.71'

# Row 22
$ws.Range("A22").Value = ' Sharjah'
$ws.Range("B22").Value = ' October 31 2020'
$ws.Range("C22").Value = 'Sunrisers won by 5 wickets (with 35 balls remaining)'
$ws.Range("D22").Value = 'Royal Challengers Bangalore'
$ws.Range("E22").Value = 'Sunrisers Hyderabad'
$ws.Range("F22").Value = 'Virat Kohli (c)'
$ws.Range("G22").Value = '7'
$ws.Range("H22").Value = '7'
$ws.Range("I22").Value = '0'
$ws.Range("J22").Value = '0'
$ws.Range("K22").Value = '100.00'

# Row 23
$ws.Range("A23").Value = ' Sharjah'
$ws.Range("B23").Value = ' October 15 2020'
$ws.Range("C23").Value = 'Kings XI won by 8 wickets'
$ws.Range("D23").Value = 'Royal Challengers Bangalore'
$ws.Range("E23").Value = 'Kings XI Punjab'
$ws.Range("F23").Value = 'Virat Kohli (c)'
$ws.Range("G23").Value = '48'
$ws.Range("H23").Value = '39'
$ws.Range("I23").Value = '3'
$ws.Range("J23").Value = '0'
$ws.Range("K23").Value = '123.07'

# Row 24
$ws.Range("A24").Value = ' Dubai (DSC)'
$ws.Range("B24").Value = ' October 05 2020'
$ws.Range("C24").Value = 'Capitals won by 59 runs'
$ws.Range("D24").Value = 'Royal Challengers Bangalore'
$ws.Range("E24").Value = 'Delhi Capitals'
$ws.Range("F24").Value = 'Virat Kohli (c)'
$ws.Range("G24").Value = '43'
$ws.Range("H24").Value = '39'
$ws.Range("I24").Value = '2'
$ws.Range("J24").Value = '1'
$ws.Range("K24").Value = '110.25'

# Row 25
$ws.Range("A25").Value = ' Dubai (DSC)'
$ws.Range("B25").Value = ' September 28 2020'
$ws.Range("C25").Value = 'Match tied (RCB won the one-over eliminator)'
$ws.Range("D25").Value = 'Royal Challengers Bangalore'
$ws.Range("E25").Value = 'Mumbai Indians'
$ws.Range("F25").Value = 'Virat Kohli (c)'
$ws.Range("G25").Value = '3'
$ws.Range("H25").Value = '11'
$ws.Range("I25").Value = '0'
$ws.Range("J25").Value = '0'
$ws.Range("K25").Value = '27.27'

# Row 26
$ws.Range("A26").Value = ' Dubai (DSC)'
$ws.Range("B26").Value = ' September 21 2020'
$ws.Range("C26").Value = 'RCB won by 10 runs'
$ws.Range("D26").Value = 'Royal Challengers Bangalore'
$ws.Range("E26").Value = 'Sunrisers Hyderabad'
$ws.Range("F26").Value = 'Virat Kohli (c)'
$ws.Range("G26").Value = '14'
$ws.Range("H26").Value = '13'
$ws.Range("I26").Value = '0'
$ws.Range("J26").Value = '0'
$ws.Range("K26").Value = '107.69'

# Row 27
$ws.Range("A27").Value = ' Dubai (DSC)'
$ws.Range("B27").Value = ' October 17 2020'
$ws.Range("C27").Value = 'RCB won by 7 wickets (with 2 balls remaining)'
$ws.Range("D27").Value = 'Royal Challengers Bangalore'
$ws.Range("E27").Value = 'Rajasthan Royals'
$ws.Range("F27").Value = 'Virat Kohli (c)'
$ws.Range("G27").Value = '43'
$ws.Range("H27").Value = '32'
$ws.Range("I27").Value = '1'
$ws.Range("J27").Value = '2'
$ws.Range("K27").Value = '134.37'

# Row 28
$ws.Range("A28").Value = ' Abu Dhabi'
$ws.Range("B28").Value = ' October 21 2020'
$ws.Range("C28").Value = 'RCB won by 8 wickets (with 39 balls remaining)'
$ws.Range("D28").Value = 'Royal Challengers Bangalore'
$ws.Range("E28").Value = 'Kolkata Knight Riders'
$ws.Range("F28").Value = 'Virat Kohli (c)'
$ws.Range("G28").Value = '18'
$ws.Range("H28").Value = '17'
$ws.Range("I28").Value = '2'
$ws.Range("J28").Value = '0'
$ws.Range("K28").Value = '105.88'

# Row 29
$ws.Range("A29").Value = ' Sharjah'
$ws.Range("B29").Value = ' October 12 2020'
$ws.Range("C29").Value = 'RCB won by 82 runs'
$ws.Range("D29").Value = 'Royal Challengers Bangalore'
$ws.Range("E29").Value = 'Kolkata Knight Riders'
$ws.Range("F29").Value = 'Virat Kohli (c)'
$ws.Range("G29").Value = '33'
$ws.Range("H29").Value = '28'
$ws.Range("I29").Value = '1'
$ws.Range("J29").Value = '0'
$ws.Range("K29").Value = '117.85'

# Row 30
$ws.Range("A30").Value = ' Dubai (DSC)'
$ws.Range("B30").Value = ' October 10 2020'
$ws.Range("C30").Value = 'RCB won by 37 runs'
$ws.Range("D30").Value = 'Royal Challengers Bangalore'
$ws.Range("E30").Value = 'Chennai Super Kings'
$ws.Range("F30").Value = 'Virat Kohli (c)'
$ws.Range("G30").Value = '90'
$ws.Range("H30").Value = '52'
$ws.Range("I30").Value = '4'
$ws.Range("J30").Value = '4'
$ws.Range("K30").Value = '173.07'

# Row 31
$ws.Range("A31").Value = ' Abu Dhabi'
$ws.Range("B31").Value = ' October 03 2020'
$ws.Range("C31").Value = 'RCB won by 8 wickets (with 5 balls remaining)'
$ws.Range("D31").Value = 'Royal Challengers Bangalore'
$ws.Range("E31").Value = 'Rajasthan Royals'
$ws.Range("F31").Value = 'Virat Kohli (c)'
$ws.Range("G31").Value = '72'
$ws.Range("H31").Value = '53'
$ws.Range("I31").Value = '7'
$ws.Range("J31").Value = '2'
$ws.Range("K31").Value = '135.84'
